$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Primer-Probe Interactions")

# --- Row 6: replace the pasted probe/primer sequence text -------------------
$ws1.Range("D6").Value = "AAGCAGAGGATGAATCT"
$ws1.Range("E6").Value = "GAAGAAGCCCCGGTCATCA"
$ws1.Range("F6").Value = "GCTTTCTTGCACGCTGGAA"
$ws1.Range("G6").Value = "rat chr1 reference"

# --- Row 7: label only changes ----------------------------------------------
$ws1.Range("G7").Value = "rat chrY target"

# --- Row 8: label only changes ----------------------------------------------
$ws1.Range("G8").Value = "rat chrX target"

# --- Row 16: new pasted sequence (same formatting as the D6:F7 paste) ------
$ws1.Range("E6").Copy($ws1.Range("F16"))
$ws1.Range("F16").Value = "GAAGAAGCCCCGGTCATCA"

# --- Column D needs to be a bit wider for the new text ----------------------
$ws1.Range("D1").ColumnWidth = 28.625

# --- Selection left where the user finished editing -------------------------
$ws1.Range("D8").Select()
